# Refresh the cryptocurrency price/volume table (columns D and E) with the
# latest scraped figures. Values that would otherwise be auto-detected as
# numbers by Excel are entered with a leading apostrophe so they remain
# plain text, matching the original inlineStr cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.681.26'
$ws.Range("E2").Value = '  -0.94%  '
$ws.Range("D3").Value = '3.787.67'
$ws.Range("E3").Value = '  -0.10%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'595.62"
$ws.Range("E5").Value = '  +0.20%  '
$ws.Range("D6").Value = "'166.68"
$ws.Range("E6").Value = '  -0.80%  '
$ws.Range("D7").Value = '3.785.15'
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  -0.23%  '
$ws.Range("E10").Value = '  -1.02%  '
$ws.Range("E11").Value = '  -0.92%  '
$ws.Range("E12").Value = '  -0.38%  '
$ws.Range("E13").Value = '  -2.97%  '
$ws.Range("E14").Value = '  -0.42%  '
$ws.Range("D15").Value = '4.423.59'
$ws.Range("E15").Value = '  -0.08%  '
$ws.Range("D16").Value = '3.791.27'
$ws.Range("E16").Value = '  -0.06%  '
$ws.Range("D17").Value = "'18.58"
$ws.Range("E17").Value = '  +3.31%  '
$ws.Range("D18").Value = '67.631.90'
$ws.Range("E18").Value = '  -1.08%  '
$ws.Range("D19").Value = "'7.05"
$ws.Range("E19").Value = '  +1.14%  '
$ws.Range("E20").Value = '  +0.07%  '
$ws.Range("E21").Value = '  -8.53%  '
$ws.Range("D22").Value = "'459.06"
$ws.Range("E22").Value = '  -1.43%  '
$ws.Range("E23").Value = '  -0.26%  '
$ws.Range("E24").Value = '  +1.81%  '
$ws.Range("D25").Value = "'83.35"
$ws.Range("E25").Value = '  -0.75%  '
$ws.Range("E26").Value = '  +0.79%  '
$ws.Range("E27").Value = '  -3.25%  '
$ws.Range("E28").Value = '  -0.09%  '
$ws.Range("D29").Value = "'10.01"
$ws.Range("E29").Value = '  -1.52%  '
$ws.Range("D30").Value = '3.935.39'
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("E31").Value = '  +0.00%  '
$ws.Range("D32").Value = "'2.27"
$ws.Range("E32").Value = '  +4.86%  '
$ws.Range("D33").Value = "'7.22"
$ws.Range("E33").Value = '  -1.44%  '
$ws.Range("D34").Value = "'29.57"
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = '  -0.09%  '
$ws.Range("E36").Value = '  -1.04%  '
$ws.Range("D37").Value = "'0.0998"
$ws.Range("E37").Value = '  -0.99%  '
$ws.Range("D38").Value = "'3.33"
$ws.Range("E38").Value = '  -2.97%  '
$ws.Range("E39").Value = '  -0.75%  '
$ws.Range("D40").Value = "'0.994"
$ws.Range("E40").Value = '  -0.98%  '
$ws.Range("E41").Value = '  -0.46%  '
$ws.Range("E42").Value = '  -0.10%  '
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("E44").Value = '  +1.98%  '
$ws.Range("D45").Value = "'43.82"
$ws.Range("E45").Value = '  -0.95%  '
$ws.Range("E46").Value = '  -1.62%  '
$ws.Range("D47").Value = "'150.57"
$ws.Range("E47").Value = '  +2.73%  '
$ws.Range("D48").Value = "'8.28"
$ws.Range("D49").Value = "'26.80"
$ws.Range("E49").Value = '  +5.55%  '
$ws.Range("D50").Value = "'388.47"
$ws.Range("E50").Value = '  -1.49%  '
$ws.Range("E51").Value = '  -5.05%  '
